$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Pass" test run status result from I2 (final execution result cleared)
$ws.Range("I2").Value = $null

# Select I1 as the active cell (matches the final saved selection in the sheet view)
$ws.Range("I1").Select()
